$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22 - this shifts the existing rows 22..100 down
# to 23..101 (matching the unified diff, which shows every subsequent row's
# data shifted down by one row).
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly price record.
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C22").Value = "Arica y Parinacota"
$ws.Range("D22").Value = 44811
$ws.Range("E22").Value = 15
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100102
$ws.Range("H22").Value = "Cítricos"
$ws.Range("I22").Value = 100102005
$ws.Range("J22").Value = "Naranja"
$ws.Range("K22").Value = "Lane Late"
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 600
$ws.Range("O22").Value = 650
$ws.Range("P22").Value = 625
$ws.Range("Q22").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R22").Value = "Región de Coquimbo"
$ws.Range("S22").Value = 625
$ws.Range("T22").Value = 1

# Match the date cell formatting used by the rest of column D.
$ws.Range("D22").NumberFormat = $ws.Range("D23").NumberFormat
